$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 22 ended at 21:00 (0.875) and is now corrected to end at 22:00.
$ws.Range("C22").Value = 22/24
